$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new worksheet "mads_efast" after "mads_tightened"
# ------------------------------------------------------------------
$tightened = $wb.Worksheets.Item("mads_tightened")
$newSheet = $wb.Worksheets.Add($null, $tightened)
$newSheet.Name = "mads_efast"

# ------------------------------------------------------------------
# 2. Fill in the formulas that mirror mads_tightened.
#    Rows 1-11 map 1:1 onto mads_tightened rows 1-11; row 12 on the
#    new sheet pulls from row 13 of mads_tightened (row 12 there,
#    "d", is intentionally skipped).
# ------------------------------------------------------------------
$cols = @("A", "B", "C", "D")

for ($r = 1; $r -le 11; $r++) {
    foreach ($col in $cols) {
        $cell = "$col$r"
        $newSheet.Range($cell).Formula = "=mads_tightened!$cell"
    }
}

foreach ($col in $cols) {
    $newSheet.Range("${col}12").Formula = "=mads_tightened!${col}13"
}

# ------------------------------------------------------------------
# 3. Sheet view / formatting to match mads_tightened's style
# ------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 20.95
$newSheet.Columns.Item(2).ColumnWidth = 10.65

$newSheet.Range("B32").Select()
$excel.ActiveWindow.Zoom = 160

# Make the new sheet the active tab of the workbook
$newSheet.Activate()
